# spring 24 week 3 inputs
# Update the SL matchup averages with the latest weekly numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nine")

$ws.Range("G2").Value = 13.5

$ws.Range("D3").Value = 10.33
$ws.Range("F3").Value = 10.26
$ws.Range("G3").Value = 11.54

$ws.Range("C4").Value = 9.67
$ws.Range("F4").Value = 10.17

$ws.Range("C6").Value = 9.74
$ws.Range("D6").Value = 9.83
$ws.Range("G6").Value = 10.33

$ws.Range("B7").Value = 6.5
$ws.Range("C7").Value = 8.46
$ws.Range("F7").Value = 9.67

$ws.Range("I8").Value = 7.62

$ws.Range("H9").Value = 12.38
